$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") updates: values like "1.000" or "0.5165" look numeric to
# Excel, so force the cell to Text format before assigning, then restore the
# cell style afterwards so no stray formatting is left behind.
$priceUpdates = @(
    @{Row=2; Value='26.397.82'}
    @{Row=3; Value='1.836.40'}
    @{Row=4; Value='1.000'}
    @{Row=5; Value='260.16'}
    @{Row=6; Value='1.000'}
    @{Row=7; Value='0.5165'}
    @{Row=8; Value='0.3220'}
    @{Row=9; Value='0.06757'}
    @{Row=10; Value='18.74'}
    @{Row=11; Value='0.7679'}
    @{Row=12; Value='0.07665'}
    @{Row=13; Value='1.840.45'}
    @{Row=14; Value='89.02'}
    @{Row=15; Value='5.042'}
    @{Row=17; Value='14.03'}
    @{Row=18; Value='0.9998'}
    @{Row=19; Value='0.000007910'}
    @{Row=20; Value='26.414.78'}
    @{Row=21; Value='2.062.26'}
    @{Row=22; Value='4.574'}
    @{Row=23; Value='9.427'}
    @{Row=24; Value='5.934'}
    @{Row=25; Value='2.290'}
    @{Row=26; Value='144.57'}
    @{Row=27; Value='1.633'}
    @{Row=29; Value='111.53'}
    @{Row=30; Value='4.198'}
    @{Row=31; Value='4.124'}
    @{Row=32; Value='0.08757'}
    @{Row=33; Value='0.04829'}
    @{Row=34; Value='1.131'}
    @{Row=35; Value='2.852'}
    @{Row=36; Value='0.6901'}
    @{Row=37; Value='3.088'}
    @{Row=38; Value='0.01775'}
    @{Row=39; Value='2.201'}
    @{Row=40; Value='0.4917'}
    @{Row=41; Value='111.94'}
    @{Row=42; Value='0.8904'}
    @{Row=43; Value='6.135'}
    @{Row=44; Value='1.0000'}
    @{Row=45; Value='7.752'}
    @{Row=46; Value='0.4193'}
    @{Row=47; Value='9.088'}
    @{Row=48; Value='0.1253'}
    @{Row=49; Value='0.05864'}
    @{Row=50; Value='35.12'}
)

$volumeUpdates = @(
    @{Row=2; Value='  -2.91%  '}
    @{Row=3; Value='  -2.58%  '}
    @{Row=4; Value='  -0.09%  '}
    @{Row=5; Value='  -7.32%  '}
    @{Row=6; Value='  -0.04%  '}
    @{Row=7; Value='  -2.29%  '}
    @{Row=8; Value='  -7.89%  '}
    @{Row=9; Value='  -3.24%  '}
    @{Row=10; Value='  -7.93%  '}
    @{Row=11; Value='  -5.89%  '}
    @{Row=12; Value='  -2.70%  '}
    @{Row=13; Value='  -0.85%  '}
    @{Row=14; Value='  -1.63%  '}
    @{Row=15; Value='  -2.53%  '}
    @{Row=16; Value='  -0.08%  '}
    @{Row=17; Value='  -3.85%  '}
    @{Row=18; Value='  -0.04%  '}
    @{Row=19; Value='  -2.81%  '}
    @{Row=20; Value='  -3.00%  '}
    @{Row=21; Value='  -3.00%  '}
    @{Row=22; Value='  -3.82%  '}
    @{Row=23; Value='  -6.70%  '}
    @{Row=24; Value='  -4.87%  '}
    @{Row=25; Value='  -3.89%  '}
    @{Row=26; Value='  -2.56%  '}
    @{Row=27; Value='  -2.37%  '}
    @{Row=28; Value='  -3.39%  '}
    @{Row=29; Value='  -4.51%  '}
    @{Row=30; Value='  -4.15%  '}
    @{Row=31; Value='  -5.27%  '}
    @{Row=32; Value='  -2.20%  '}
    @{Row=33; Value='  -2.45%  '}
    @{Row=34; Value='  -4.54%  '}
    @{Row=35; Value='  -1.74%  '}
    @{Row=36; Value='  -6.96%  '}
    @{Row=37; Value='  -7.31%  '}
    @{Row=38; Value='  -5.18%  '}
    @{Row=39; Value='  -8.35%  '}
    @{Row=40; Value='  -6.24%  '}
    @{Row=41; Value='  -3.37%  '}
    @{Row=42; Value='  -8.34%  '}
    @{Row=43; Value='  -1.78%  '}
    @{Row=44; Value='  +0.01%  '}
    @{Row=45; Value='  -4.56%  '}
    @{Row=46; Value='  -7.92%  '}
    @{Row=47; Value='  -3.77%  '}
    @{Row=48; Value='  -7.82%  '}
    @{Row=49; Value='  -1.63%  '}
    @{Row=50; Value='  -3.97%  '}
)

foreach ($update in $priceUpdates) {
    $cell = $ws.Cells.Item($update.Row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $update.Value
    $cell.Style = "Normal"
}

foreach ($update in $volumeUpdates) {
    $ws.Cells.Item($update.Row, 5).Value = $update.Value
}

# Row 51: the coin listed dropped out and was replaced by a new entry (Aave -> EOS)
$ws.Range("B51").Value = "EOS"
$ws.Range("C51").Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.8818"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.69%  "
